$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Distributed System & Database*") {
        $full = $p.Range
        $xmlFrag = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="07316229" w14:textId="562B079D" w:rsidR="0079475D" w:rsidRPr="00470965" w:rsidRDefault="00D07A66" w:rsidP="00FF341E"><w:pPr><w:spacing w:before="60" w:line="260" w:lineRule="exact"/><w:jc w:val="left"/><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r w:rsidRPr="008F7912"><w:rPr><w:rFonts w:hint="eastAsia"/><w:b/><w:bCs/><w:szCs w:val="21"/></w:rPr><w:t>Distributed System &amp; Database</w:t></w:r><w:r w:rsidRPr="008F7912"><w:rPr><w:rFonts w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00026CAF" w:rsidRPr="00026CAF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">Developed a </w:t></w:r><w:r w:rsidR="00DF0C39"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri" w:hint="eastAsia"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="21"/></w:rPr><w:t>distributed</w:t></w:r><w:r w:rsidR="00026CAF" w:rsidRPr="00026CAF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">, sharded key/value store in Java using Multi-Paxos, later optimizing it with a Raft-like consensus algorithm for improved </w:t></w:r><w:r w:rsidR="001520AB" w:rsidRPr="001520AB"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="21"/></w:rPr><w:t>throughput</w:t></w:r><w:r w:rsidR="00026CAF" w:rsidRPr="00026CAF"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="21"/></w:rPr><w:t>; implemented a two-phase commit protocol to ensure atomic, strongly consistent cross-shard transactions in a horizontally scalable system.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri" w:hint="eastAsia"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> [</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri" w:hint="eastAsia"/><w:szCs w:val="21"/></w:rPr><w:t>Link</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri" w:hint="eastAsia"/><w:color w:val="000000" w:themeColor="text1"/><w:szCs w:val="21"/></w:rPr><w:t>]</w:t></w:r></w:p>
"@
        $full.InsertXML($xmlFrag)
        break
    }
}

# Locate the newly-inserted "Link" run precisely using the unique surrounding context.
$search = $d.Content
$search.Find.ClearFormatting()
$search.Find.Execute("horizontally scalable system. [Link]")
$matchStart = $search.Start
$matchText = $search.Text
$linkOffset = $matchText.IndexOf("[Link]") + 1
$linkStart = $matchStart + $linkOffset
$linkEnd = $linkStart + 4
$linkRng = $d.Range($linkStart, $linkEnd)

$d.Hyperlinks.Add($linkRng, "https://github.com/yimianxyz/dslab", [Type]::Missing, [Type]::Missing, "Link") | Out-Null

Write-Host "done"
